# Auto-generated edit script applying numeric corrections to Sheets/Lich_Profits.xlsx
# (workbook sheet tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 20004444
$ws.Range("I40").Value = 6476.857
$ws.Range("J40").Value = 37502664
$ws.Range("K40").Value = 6476.857
$ws.Range("L40").Value = 37502664
$ws.Range("M40").Value = -6301.857
$ws.Range("N40").Value = -37503014

# Row 86
$ws.Range("H86").Value = 6378.9
$ws.Range("I86").Value = 6421
$ws.Range("K86").Value = 6421
$ws.Range("M86").Value = -5298

# Row 89
$ws.Range("H89").Value = 6378.9
$ws.Range("I89").Value = 6421
$ws.Range("K89").Value = 32105
$ws.Range("M89").Value = -26489

# Row 113
$ws.Range("H113").Value = 5553.6924
$ws.Range("I113").Value = 5824.6665
$ws.Range("J113").Value = 4944
$ws.Range("K113").Value = 5824.6665
$ws.Range("L113").Value = 4944
$ws.Range("M113").Value = -2570.6665
$ws.Range("N113").Value = -11452

# Row 141
$ws.Range("H141").Value = 2556.25
$ws.Range("I141").Value = 2660.3333
$ws.Range("J141").Value = 995
$ws.Range("K141").Value = 7980.999899999999
$ws.Range("L141").Value = 2985
$ws.Range("M141").Value = -2800.999899999999
$ws.Range("N141").Value = -13345

$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 5251.3335
$ws.Range("I11").Value = 103
$ws.Range("J11").Value = 6281
$ws.Range("K11").Value = 103
$ws.Range("L11").Value = 6281
$ws.Range("M11").Value = 41
$ws.Range("N11").Value = -6569

# Row 13
$ws.Range("H13").Value = 1670.8334
$ws.Range("J13").Value = 405
$ws.Range("L13").Value = 405
$ws.Range("N13").Value = -693

# Row 45
$ws.Range("H45").Value = 2039.5
$ws.Range("I45").Value = 1993.5
$ws.Range("K45").Value = 1993.5
$ws.Range("M45").Value = -1616.5

# Row 61
$ws.Range("H61").Value = 2855.375
$ws.Range("I61").Value = 2120.1333
$ws.Range("K61").Value = 2120.1333
$ws.Range("M61").Value = -1908.1333

# Row 102
$ws.Range("H102").Value = 1433.381
$ws.Range("I102").Value = 1294.8422
$ws.Range("K102").Value = 1294.8422
$ws.Range("M102").Value = 327.1578

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

# Row 122
$ws.Range("H122").Value = 2765.2222
$ws.Range("I122").Value = 2657.5945
$ws.Range("J122").Value = 3263
$ws.Range("K122").Value = 7972.7835
$ws.Range("L122").Value = 9789
$ws.Range("M122").Value = -5522.7835
$ws.Range("N122").Value = -14689

# Row 132
$ws.Range("H132").Value = 2976.6924
$ws.Range("I132").Value = 2976.6924
$ws.Range("K132").Value = 8930.0772
$ws.Range("M132").Value = -6400.0772

# Row 136
$ws.Range("H136").Value = 2855.375
$ws.Range("I136").Value = 2120.1333
$ws.Range("K136").Value = 6360.3999
$ws.Range("M136").Value = -3810.3999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2495.6667
$ws.Range("I86").Value = 2495.6667
$ws.Range("K86").Value = 2495.6667
$ws.Range("M86").Value = -1372.6667

# Row 89
$ws.Range("H89").Value = 2495.6667
$ws.Range("I89").Value = 2495.6667
$ws.Range("K89").Value = 12478.3335
$ws.Range("M89").Value = -6862.333500000001

# Row 134
$ws.Range("H134").Value = 3386.356
$ws.Range("I134").Value = 3855.1667
$ws.Range("K134").Value = 11565.5001
$ws.Range("M134").Value = -9030.500100000001

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 407.63635
$ws.Range("I22").Value = 421
$ws.Range("J22").Value = 347.5
$ws.Range("K22").Value = 421
$ws.Range("L22").Value = 347.5
$ws.Range("M22").Value = -71
$ws.Range("N22").Value = -1047.5

# Row 41
$ws.Range("H41").Value = 7291.8
$ws.Range("I41").Value = 1029.5
$ws.Range("J41").Value = 11466.667
$ws.Range("K41").Value = 1029.5
$ws.Range("L41").Value = 11466.667
$ws.Range("M41").Value = -601.5
$ws.Range("N41").Value = -12322.667

# Row 47
$ws.Range("H47").Value = 8500.5
$ws.Range("I47").Value = 8500.5
$ws.Range("K47").Value = 8500.5
$ws.Range("M47").Value = -7934.5

# Row 122
$ws.Range("H122").Value = 1897
$ws.Range("J122").Value = 1895
$ws.Range("L122").Value = 5685
$ws.Range("N122").Value = -10585

# Row 132
$ws.Range("H132").Value = 4678.9663
$ws.Range("I132").Value = 4596.2173
$ws.Range("K132").Value = 13788.6519
$ws.Range("M132").Value = -11258.6519

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 84.954544
$ws.Range("I40").Value = 92.94118
$ws.Range("J40").Value = 57.8
$ws.Range("K40").Value = 371.76472
$ws.Range("L40").Value = 231.2
$ws.Range("M40").Value = -302.76472
$ws.Range("N40").Value = -369.2

# Row 68
$ws.Range("H68").Value = 3335120
$ws.Range("J68").Value = 2223.2778
$ws.Range("L68").Value = 6669.8334
$ws.Range("N68").Value = -8291.8334

# Row 71
$ws.Range("H71").Value = 3335120
$ws.Range("J71").Value = 2223.2778
$ws.Range("L71").Value = 20009.5002
$ws.Range("N71").Value = -28121.5002

# Row 107
$ws.Range("H107").Value = 1061.4
$ws.Range("I107").Value = 542.875
$ws.Range("J107").Value = 1249.9546
$ws.Range("K107").Value = 1628.625
$ws.Range("L107").Value = 3749.8638
$ws.Range("M107").Value = 291.375
$ws.Range("N107").Value = -7589.8638

# Row 140
$ws.Range("H140").Value = 10770
$ws.Range("I140").Value = 12657.529
$ws.Range("J140").Value = 5422
$ws.Range("K140").Value = 37972.587
$ws.Range("L140").Value = 16266
$ws.Range("M140").Value = -32792.587
$ws.Range("N140").Value = -26626

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10304.333
$ws.Range("I70").Value = 7992
$ws.Range("J70").Value = 13194.75
$ws.Range("K70").Value = 7992
$ws.Range("L70").Value = 13194.75
$ws.Range("M70").Value = -7722
$ws.Range("N70").Value = -13734.75

# Row 73
$ws.Range("H73").Value = 10304.333
$ws.Range("I73").Value = 7992
$ws.Range("J73").Value = 13194.75
$ws.Range("K73").Value = 7992
$ws.Range("L73").Value = 13194.75
$ws.Range("M73").Value = -7056
$ws.Range("N73").Value = -15066.75

# Row 80
$ws.Range("H80").Value = 2369.4285
$ws.Range("J80").Value = 2649
$ws.Range("L80").Value = 2649
$ws.Range("N80").Value = -4645

# Row 83
$ws.Range("H83").Value = 2369.4285
$ws.Range("J83").Value = 2649
$ws.Range("L83").Value = 13245
$ws.Range("N83").Value = -23229

# Row 102
$ws.Range("H102").Value = 3423.9644
$ws.Range("J102").Value = 5351.2
$ws.Range("L102").Value = 5351.2
$ws.Range("N102").Value = -8595.200000000001

# Row 122
$ws.Range("H122").Value = 2965.625
$ws.Range("I122").Value = 3023.5454
$ws.Range("K122").Value = 9070.636200000001
$ws.Range("M122").Value = -6620.636200000001

# Row 126
$ws.Range("H126").Value = 10614.667
$ws.Range("I126").Value = 40000
$ws.Range("J126").Value = 4737.6
$ws.Range("K126").Value = 120000
$ws.Range("L126").Value = 14212.8
$ws.Range("M126").Value = -117530
$ws.Range("N126").Value = -19152.8

# Row 132
$ws.Range("H132").Value = 61764.176
$ws.Range("I132").Value = 79600.53999999999
$ws.Range("J132").Value = 3796
$ws.Range("K132").Value = 238801.62
$ws.Range("L132").Value = 11388
$ws.Range("M132").Value = -236271.62
$ws.Range("N132").Value = -16448

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 6974.75
$ws.Range("I3").Value = 4500
$ws.Range("K3").Value = 4500
$ws.Range("M3").Value = -4388

# Row 15
$ws.Range("H15").Value = 6974.75
$ws.Range("I15").Value = 4500
$ws.Range("K15").Value = 4500
$ws.Range("M15").Value = -4330

# Row 46
$ws.Range("H46").Value = 3397.6
$ws.Range("J46").Value = 4898
$ws.Range("L46").Value = 4898
$ws.Range("N46").Value = -5274

# Row 55
$ws.Range("H55").Value = 599.43335
$ws.Range("I55").Value = 721.6842
$ws.Range("J55").Value = 388.27274
$ws.Range("K55").Value = 721.6842
$ws.Range("L55").Value = 388.27274
$ws.Range("M55").Value = -548.6842
$ws.Range("N55").Value = -734.27274

# Row 93
$ws.Range("H93").Value = 1758
$ws.Range("I93").Value = 1455
$ws.Range("K93").Value = 1455
$ws.Range("M93").Value = -207

# Row 122
$ws.Range("H122").Value = 5266.5
$ws.Range("I122").Value = 5540
$ws.Range("K122").Value = 16620
$ws.Range("M122").Value = -14170

# Row 132
$ws.Range("H132").Value = 4559.846
$ws.Range("I132").Value = 4470
$ws.Range("J132").Value = 4859.3335
$ws.Range("K132").Value = 13410
$ws.Range("L132").Value = 14578.0005
$ws.Range("M132").Value = -10880
$ws.Range("N132").Value = -19638.0005

# Row 136
$ws.Range("H136").Value = 3181.3333
$ws.Range("I136").Value = 2992.8572
$ws.Range("J136").Value = 3384.3076
$ws.Range("K136").Value = 8978.571599999999
$ws.Range("L136").Value = 10152.9228
$ws.Range("M136").Value = -6428.571599999999
$ws.Range("N136").Value = -15252.9228

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 4500
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = ""

# Row 14
$ws.Range("H14").Value = 22249.75
$ws.Range("J14").Value = 17999.666
$ws.Range("L14").Value = 17999.666
$ws.Range("N14").Value = -18335.666

# Row 122
$ws.Range("H122").Value = 6863.231
$ws.Range("I122").Value = 7205.364
$ws.Range("J122").Value = 4981.5
$ws.Range("K122").Value = 21616.092
$ws.Range("L122").Value = 14944.5
$ws.Range("M122").Value = -19166.092
$ws.Range("N122").Value = -19844.5

# Row 132
$ws.Range("H132").Value = 1812.4559
$ws.Range("I132").Value = 1228.1464
$ws.Range("J132").Value = 2699.7407
$ws.Range("K132").Value = 3684.4392
$ws.Range("L132").Value = 8099.222099999999
$ws.Range("M132").Value = -1154.4392
$ws.Range("N132").Value = -13159.2221

# Row 136
$ws.Range("H136").Value = 335612.3
$ws.Range("I136").Value = 347150.7
$ws.Range("J136").Value = 999
$ws.Range("K136").Value = 1041452.1
$ws.Range("L136").Value = 2997
$ws.Range("M136").Value = -1038902.1
$ws.Range("N136").Value = -8097
